$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.523
$ws.Range("E4").Value = 13.265
$ws.Range("A9").Value = -20.912
$ws.Range("E10").Value = 12.554
$ws.Range("A18").Value = -21.81
$ws.Range("A20").Value = -21.828
$ws.Range("D21").Value = -7.675999999999999
